# Weekly update: insert a new Coliflor price record for Vega Modelo de Temuco
# as row 197 (Fecha 2021-10-07 / 44476, Volumen 3500, Origen "Región Metropolitana"),
# pushing the existing rows 197-239 down to 198-240.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 197 - shifts rows 197:239 down to 198:240
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row with the new observation
$ws.Cells.Item(197, 1).Value = 10
$ws.Cells.Item(197, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(197, 3).Value = "La Araucanía"
$ws.Cells.Item(197, 4).Value = 44476
$ws.Cells.Item(197, 5).Value = 9
$ws.Cells.Item(197, 6).Value = 100112008
$ws.Cells.Item(197, 7).Value = "Coliflor"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 3500
$ws.Cells.Item(197, 11).Value = 800
$ws.Cells.Item(197, 12).Value = 800
$ws.Cells.Item(197, 13).Value = 800
$ws.Cells.Item(197, 14).Value = "$/unidad"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 800
$ws.Cells.Item(197, 17).Value = 1
$ws.Cells.Item(197, 18).Value = "Hortaliza"
